$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Fix capitalization of "de"/"del" -> "De"/"Del" in place names ---
$ws.Range("B2").Value = "Comitán De Domínguez"
$ws.Range("B5").Value = "Hidalgo Del Parral"
$ws.Range("A13").Value = "Ciudad De México"
$ws.Range("A24").Value = "Estado De México"
$ws.Range("B24").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B26").Value = "Ecatepec De Morelos"
$ws.Range("B32").Value = "Acapulco De Juárez"
$ws.Range("B35").Value = "Mártir De Cuilapan"
$ws.Range("B36").Value = "Tixtla De Guerrero"
$ws.Range("B42").Value = "Autlán De Navarro"
$ws.Range("B44").Value = "Cuautitlán De García Barragán"
$ws.Range("B63").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B75").Value = "Cadereyta De Montes"
$ws.Range("B77").Value = "Landa De Matamoros"
$ws.Range("B86").Value = "Ixhuatlán Del Sureste"
$ws.Range("B87").Value = "Soledad De Doblado"

# --- Remove trailing footer/metadata rows no longer needed ---
# Rows 97-101 (footer block right after the first data section)
$ws.Rows("97:101").Delete()
# Rows 476-480 (trailing footer block) shift up to 471-475 once the rows above are removed
$ws.Rows("471:475").Delete()
